$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tests")

# Widen column A to fit the new, longer test case text
$ws.Columns.Item(1).ColumnWidth = 97.5

# New rows of data appended after row 161 (row 162 stays blank, like the
# blank separator rows before every other "Test:" block in the sheet)
$data = @(
    @(163, 'A', 'Test: Bildschirme überspringen'),
    @(164, 'A', 'Der Nutzer wählt für Spieler 1 "Mensch", für Spieler 2 "KI-Elimination" und klickt auf das Feld "Spiel starten".'),
    @(164, 'B', 'Die Spielansicht wird angezeigt.'),
    @(165, 'A', 'Der Nutzer klickt auf die Checkbox "Startansicht überspringen".'),
    @(165, 'B', 'Die Checkbox wird mit einem Haken versehen.'),
    @(166, 'A', 'Der Nutzer klickt abwechselnd selber auf ein freies Feld und löst Züge der KI aus, bis das Spiel beendet ist.'),
    @(167, 'A', 'Der Nutzer klickt auf das Feld "Weiter".'),
    @(167, 'B', 'Es wird in die Belohnungsansicht gewechselt.'),
    @(168, 'A', 'Der Nutzer klickt auf das Feld "Belohnung ausführen".'),
    @(168, 'B', 'Die Gewichte des Graphen ändern sich, wenn kein Unentschieden stattgefunden hat.'),
    @(169, 'A', 'Der Nutzer klickt auf das Feld "Weiter".'),
    @(169, 'B', 'Die Spielansicht wird angezeigt.'),
    @(170, 'A', 'Der Nutzer klickt auf die Checkbox "Startansicht überspringen".'),
    @(170, 'B', 'Die Checkbox ist nicht mehr mit einem Haken versehen.'),
    @(171, 'A', ' Der Nutzer klickt auf die Checkbox "Belohnungsansicht überspringen".'),
    @(171, 'B', 'Die Checkbox wird mit einem Haken versehen.'),
    @(172, 'A', 'Der Nutzer klickt abwechselnd selber auf ein freies Feld und löst Züge der KI aus, bis das Spiel beendet ist.'),
    @(173, 'A', 'Der Nutzer klickt auf das Feld "Weiter".'),
    @(173, 'B', 'Die Startansicht wird angezeigt'),
    @(174, 'A', 'Der Nutzer klickt auf das Feld "Spiel starten".'),
    @(174, 'B', 'Die Spielansicht wird angezeigt.'),
    @(175, 'A', ' Der Nutzer klickt auf die Checkbox "Belohnungsansicht überspringen".'),
    @(175, 'B', 'Die Checkbox ist nicht mehr mit einem Haken versehen.'),
    @(176, 'A', 'Der Nutzer klickt abwechselnd selber auf ein freies Feld und löst Züge der KI aus, bis das Spiel beendet ist.'),
    @(177, 'A', 'Der Nutzer klickt auf das Feld "Weiter".'),
    @(177, 'B', 'Es wird in die Belohnungsansicht gewechselt.'),
    @(178, 'A', 'Der Nutzer klickt auf das Feld "Belohnung ausführen".'),
    @(178, 'B', 'Die Gewichte des Graphen ändern sich, wenn kein Unentschieden stattgefunden hat.'),
    @(179, 'A', 'Der Nutzer klickt auf das Feld "Weiter".'),
    @(179, 'B', 'Die Startansicht wird angezeigt')
)

foreach ($entry in $data) {
    $row = $entry[0]
    $col = $entry[1]
    $text = $entry[2]
    $ws.Range("$col$row").Value = $text
}

# The "Test: ..." header cell uses the bold style applied to every other
# header row in this sheet (e.g. A156 "Test: Speedrun").
$ws.Range("A163").Font.Bold = $true

# Update the view so the newly added rows are visible, matching the
# scrolled/selected state recorded after the edit.
$ws.Range("A190").Select()
